$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Indicator cell (B4) text to the new wording
$ws.Range("B4").Value = "12.b.1 Implementation of standard accounting tools to monitor the economic and environmental aspects of tourism sustainability"

# Move the active selection to B4 (as saved in the edited workbook)
$ws.Range("B4").Select()
